$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.236.52"
$ws.Range("E2").Value = "  -0.41%  "

$ws.Range("D3").Value = "1.862.38"
$ws.Range("E3").Value = "  -0.54%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "236.30"
$ws.Range("E5").Value = "  +0.32%  "

$ws.Range("E6").Value = "  +0.05%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4680"
$ws.Range("E7").Value = "  +0.39%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2889"
$ws.Range("E8").Value = "  +1.74%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06533"
$ws.Range("E9").Value = "  -0.15%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.59"
$ws.Range("E10").Value = "  +1.52%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07942"
$ws.Range("E11").Value = "  +0.22%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "97.75"
$ws.Range("E12").Value = "  +0.66%  "

$ws.Range("D13").Value = "1.867.09"
$ws.Range("E13").Value = "  -0.28%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.161"
$ws.Range("E14").Value = "  +0.25%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6803"
$ws.Range("E15").Value = "  +0.91%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "267.30"
$ws.Range("E16").Value = "  -5.11%  "

$ws.Range("D17").Value = "30.228.26"
$ws.Range("E17").Value = "  -0.46%  "

$ws.Range("E18").Value = "  +8.11%  "

$ws.Range("E19").Value = "  +0.07%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007396"
$ws.Range("E20").Value = "  +1.57%  "

$ws.Range("D21").Value = "2.111.60"
$ws.Range("E21").Value = "  -0.72%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.308"
$ws.Range("E22").Value = "  -4.26%  "

$ws.Range("E23").Value = "  +0.00%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.165"
$ws.Range("E24").Value = "  -0.54%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "167.13"
$ws.Range("E25").Value = "  +1.88%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.190"
$ws.Range("E26").Value = "  -1.10%  "

$ws.Range("E27").Value = "  -1.27%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.955"
$ws.Range("E28").Value = "  +1.13%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.387"
$ws.Range("E29").Value = "  +2.41%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09812"
$ws.Range("E30").Value = "  +1.44%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.372"
$ws.Range("E31").Value = "  -1.24%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.470"
$ws.Range("E32").Value = "  -0.41%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.043"
$ws.Range("E33").Value = "  -1.52%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04698"
$ws.Range("E34").Value = "  +0.04%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.129"
$ws.Range("E35").Value = "  +0.90%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7019"
$ws.Range("E36").Value = "  -0.27%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.710"
$ws.Range("E37").Value = "  -0.33%  "

$ws.Range("E38").Value = "  +0.78%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.617"
$ws.Range("E39").Value = "  +2.84%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.253"
$ws.Range("E40").Value = "  -1.61%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "74.24"
$ws.Range("E41").Value = "  +0.68%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.932"
$ws.Range("E42").Value = "  -0.60%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8474"
$ws.Range("E43").Value = "  -0.10%  "

$ws.Range("B44").Value = "PaxDollar"
$ws.Range("C44").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.001"
$ws.Range("E44").Value = "  +0.02%  "

$ws.Range("B45").Value = "TheSandbox"
$ws.Range("C45").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4157"
$ws.Range("E45").Value = "  -0.73%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "103.04"
$ws.Range("E46").Value = "  -0.62%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "957.42"
$ws.Range("E47").Value = "  +2.72%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.136"
$ws.Range("E48").Value = "  -1.16%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.186"
$ws.Range("E49").Value = "  -0.82%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "34.10"
$ws.Range("E50").Value = "  -0.11%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05658"
$ws.Range("E51").Value = "  +0.57%  "

Write-Host "Updated cryptos worksheet."
